$d = $word.ActiveDocument

# 1) Remove the old "@prints" run from the "PREENCHIMENTO DO DESENVOLVIMENTO"
#    table (it used to live under "Descrição ou Imagem da Solução").
$d.Content.Find.Execute("@prints", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2) Remove the "@usuario" run from the "PREENCHIMENTO DO TESTE E QUALIDADE" table.
$d.Content.Find.Execute("@usuario", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3) Move the "@prints" placeholder: add it to the (currently empty) paragraph
#    right below "Descrição ou Imagem (menu) da Localização da Funcionalidade:"
#    in the first table (row 11, col 1).
$tbl1 = $d.Tables.Item(1)
$targetCell = $tbl1.Rows.Item(11).Cells.Item(1)
$targetCell.Range.InsertAfter("@prints")
